# Update the "Anna Milligan" credit line on the "Thank you" slide so it
# reads "Anna Milligan – Frontend, Reports" (split into two runs, matching
# how the second contributor line above it already lists multiple areas).

$p = $ppt.ActivePresentation

$targetShape = $null
$targetParagraph = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($k = 1; $k -le $paraCount; $k++) {
                $para = $tr.Paragraphs($k, 1)
                if ($para.Text -like "Anna Milligan*Frontend*") {
                    $targetShape = $shape
                    $targetParagraph = $para
                }
            }
        }
    }
}

# First run currently holds the whole line, e.g. "Anna Milligan – Frontend".
$firstRun = $targetParagraph.Runs(1)

# Append the new trailing text right after the existing run, then shrink the
# original run down to just "Anna Milligan " so the dash/role text becomes
# its own run (mirrors how the line now reads "Frontend, Reports").
$null = $firstRun.InsertAfter([char]0x2013 + " Frontend, Reports")
$firstRun.Text = "Anna Milligan "
